$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: copy formatting from E1 (bold/centered/bordered header style)
# then set its text, matching the style of the other header cells.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:38:39.852210",
    "2021-10-05 13:38:39.852224",
    "2021-10-05 13:38:39.852228",
    "2021-10-05 13:38:39.852231",
    "2021-10-05 13:38:39.852234",
    "2021-10-05 13:38:39.852237",
    "2021-10-05 13:38:39.852240",
    "2021-10-05 13:38:39.852243",
    "2021-10-05 13:38:39.852246",
    "2021-10-05 13:38:39.852250",
    "2021-10-05 13:38:39.852253",
    "2021-10-05 13:38:39.852256",
    "2021-10-05 13:38:39.852258",
    "2021-10-05 13:38:39.852262",
    "2021-10-05 13:38:39.852265",
    "2021-10-05 13:38:39.852267",
    "2021-10-05 13:38:39.852271",
    "2021-10-05 13:38:39.852274",
    "2021-10-05 13:38:39.852277",
    "2021-10-05 13:38:39.852280"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
